# Mark attendance-derived summary cells that flip from 0 to 1.
# Rows 3-18 hold per-date attendance stats in columns D:H
# (Total Attendance Count, Real, Duplicate, Invalid, Absent).
# For most dates only the "Absent" count (column H) is set; for a few
# dates the "Total Attendance Count"/"Real" counts (columns D/E) are set
# instead, and for the very first date the "Invalid"/"Absent" counts
# (columns G/H) are set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 1

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1

$ws.Range("D5").Value = 1
$ws.Range("E5").Value = 1

$ws.Range("H6").Value = 1
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 1
$ws.Range("E9").Value = 1

$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 1
$ws.Range("E11").Value = 1

$ws.Range("H12").Value = 1
$ws.Range("H13").Value = 1
$ws.Range("H14").Value = 1
$ws.Range("H15").Value = 1
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
